# Update stats for 2025-10 (row 23 in Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B23").Value = 6312
$ws.Range("C23").Value = 1003
$ws.Range("D23").Value = 5912553
$ws.Range("E23").Value = 936.7162547528517
$ws.Range("F23").Value = 8.304735758407688
$ws.Range("G23").Value = 4.370447450572312
$ws.Range("H23").Value = 26.69087681728801
